$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3-7: g4_utils_v32.py group -> add "couverture de tests" column F ---
$covA = $ws.Range("F3:F7")
$covA.Merge()
$ws.Range("F3").Value = "80%`n(25 warnings`nsur 125 lignes de code)"
$covA.HorizontalAlignment = -4108
$covA.VerticalAlignment = -4108
$covA.WrapText = $true

# D7 loses its content (keeps its purple-font style)
$ws.Range("D7").ClearContents()

# --- Row 8-11: parserV1.py group updates ---
$ws.Range("C8").Value = (Get-Date -Year 2018 -Month 1 -Day 17).Date
$ws.Range("D8").Value = "les exceptions sont traitées (on affiche leur nom)"
$ws.Range("E8").Value = "bug"

$covB = $ws.Range("F8:F11")
$covB.Merge()
$ws.Range("F8").Value = "96%`n(2 warnings`net 1 lignes de code mort`nsur 70 lignes de code)"
$covB.HorizontalAlignment = -4108
$covB.VerticalAlignment = -4108
$covB.WrapText = $true

$ws.Range("D9").Value = 'le message "<nom_journal> OK"'
$ws.Range("E9").Value = "ok"

$ws.Range("D10").ClearContents()
$ws.Range("D11").ClearContents()

# --- Explicit row heights (Excel recorded these after the content/format edits) ---
$ws.Rows.Item(3).RowHeight = 15
$ws.Rows.Item(8).RowHeight = 40.5

# --- Selection moves from D12 to D11 ---
$ws.Range("D11").Select()
